# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Update the OFF sheet (row 3: Short Att, Short Comp, Deep Att, Deep Comp)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 508
$wsOff.Range("C3").Value = 365
$wsOff.Range("D3").Value = 127
$wsOff.Range("E3").Value = 69

# Update the DEF sheet (row 3: Short Att, Short Comp, Deep Att, Deep Comp)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 444
$wsDef.Range("C3").Value = 312
$wsDef.Range("D3").Value = 114
$wsDef.Range("E3").Value = 63
